# Weekly update: insert the new price record for Poroto verde
# (Terminal Hortofrutícola Agro Chillán) as a new row at position 139,
# pushing the existing rows 139-160 down to 140-161.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 139 - this shifts rows 139:160 down to 140:161
$ws.Rows("139:139").Insert()

# Populate the new row 139 with this week's new data point
$ws.Cells.Item(139, 1).Value2  = 7
$ws.Cells.Item(139, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(139, 3).Value2  = "Ñuble"
$ws.Cells.Item(139, 4).Value2  = 45142
$ws.Cells.Item(139, 5).Value2  = 16
$ws.Cells.Item(139, 6).Value2  = 100112031
$ws.Cells.Item(139, 7).Value2  = "Poroto verde"
$ws.Cells.Item(139, 8).Value2  = "Magnum"
$ws.Cells.Item(139, 9).Value2  = "Primera"
$ws.Cells.Item(139, 10).Value2 = 50
$ws.Cells.Item(139, 11).Value2 = 20000
$ws.Cells.Item(139, 12).Value2 = 20000
$ws.Cells.Item(139, 13).Value2 = 20000
$ws.Cells.Item(139, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(139, 15).Value2 = "Perú"
$ws.Cells.Item(139, 16).Value2 = 800
$ws.Cells.Item(139, 17).Value2 = 25
$ws.Cells.Item(139, 18).Value2 = "Hortaliza"
